$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.528161333333333
$ws.Range("H2").Value = 13.584484
$ws.Range("I2").Value = 0.3225352762763812
$ws.Range("J2").Value = 0.3225352762763812
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.102580333333333
$ws.Range("N2").Value = 3.307741
$ws.Range("O2").Value = 0.3703441503267132
$ws.Range("P2").Value = 0.3703441503267133
$ws.Range("Q2").Value = 4.992661632293777
$ws.Range("R2").Value = 44.933954690644
$ws.Range("S2").Value = 0.1194490528429681
$ws.Range("T2").Value = 0.1194490528429681

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.528161333333333
$ws.Range("H3").Value = 13.584484
$ws.Range("I3").Value = 0.3225352762763812
$ws.Range("J3").Value = 0.3225352762763812
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.200765
$ws.Range("N3").Value = 3.602295
$ws.Range("O3").Value = 0.4033232592881871
$ws.Range("P3").Value = 0.4033232592881872
$ws.Range("Q3").Value = 5.437257643420001
$ws.Range("R3").Value = 48.93531879078
$ws.Range("S3").Value = 0.130085978863206
$ws.Range("T3").Value = 0.130085978863206

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.528161333333333
$ws.Range("H4").Value = 13.584484
$ws.Range("I4").Value = 0.3225352762763812
$ws.Range("J4").Value = 0.3225352762763812
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.6738323333333334
$ws.Range("N4").Value = 2.021497
$ws.Range("O4").Value = 0.2263325903850996
$ws.Range("P4").Value = 0.2263325903850996
$ws.Range("Q4").Value = 3.051221516949778
$ws.Range("R4").Value = 27.460993652548
$ws.Range("S4").Value = 0.07300024457020714
$ws.Range("T4").Value = 0.07300024457020714

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.270951666666666
$ws.Range("H5").Value = 18.812855
$ws.Range("I5").Value = 0.4466720550425397
$ws.Range("J5").Value = 0.4466720550425397
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.102580333333333
$ws.Range("N5").Value = 3.307741
$ws.Range("O5").Value = 0.3703441503267132
$ws.Range("P5").Value = 0.3703441503267133
$ws.Range("Q5").Value = 6.914227978950555
$ws.Range("R5").Value = 62.22805181055499
$ws.Range("S5").Value = 0.1654223826994162
$ws.Range("T5").Value = 0.1654223826994163

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.270951666666666
$ws.Range("H6").Value = 18.812855
$ws.Range("I6").Value = 0.4466720550425397
$ws.Range("J6").Value = 0.4466720550425397
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.200765
$ws.Range("N6").Value = 3.602295
$ws.Range("O6").Value = 0.4033232592881871
$ws.Range("P6").Value = 0.4033232592881872
$ws.Range("Q6").Value = 7.529939278025
$ws.Range("R6").Value = 67.769453502225
$ws.Range("S6").Value = 0.1801532290727096
$ws.Range("T6").Value = 0.1801532290727097

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.270951666666666
$ws.Range("H7").Value = 18.812855
$ws.Range("I7").Value = 0.4466720550425397
$ws.Range("J7").Value = 0.4466720550425397
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.6738323333333334
$ws.Range("N7").Value = 2.021497
$ws.Range("O7").Value = 0.2263325903850996
$ws.Range("P7").Value = 0.2263325903850996
$ws.Range("Q7").Value = 4.225569993770556
$ws.Range("R7").Value = 38.030129943935
$ws.Range("S7").Value = 0.1010964432704138
$ws.Range("T7").Value = 0.1010964432704138

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.240161666666667
$ws.Range("H8").Value = 9.720485
$ws.Range("I8").Value = 0.2307926686810791
$ws.Range("J8").Value = 0.2307926686810791
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.102580333333333
$ws.Range("N8").Value = 3.307741
$ws.Range("O8").Value = 0.3703441503267132
$ws.Range("P8").Value = 0.3703441503267133
$ws.Range("Q8").Value = 3.572538530487222
$ws.Range("R8").Value = 32.152846774385
$ws.Range("S8").Value = 0.08547271478432887
$ws.Range("T8").Value = 0.08547271478432887

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.240161666666667
$ws.Range("H9").Value = 9.720485
$ws.Range("I9").Value = 0.2307926686810791
$ws.Range("J9").Value = 0.2307926686810791
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.200765
$ws.Range("N9").Value = 3.602295
$ws.Range("O9").Value = 0.4033232592881871
$ws.Range("P9").Value = 0.4033232592881872
$ws.Range("Q9").Value = 3.890672723675
$ws.Range("R9").Value = 35.01605451307501
$ws.Range("S9").Value = 0.09308405135227153
$ws.Range("T9").Value = 0.09308405135227153

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.240161666666667
$ws.Range("H10").Value = 9.720485
$ws.Range("I10").Value = 0.2307926686810791
$ws.Range("J10").Value = 0.2307926686810791
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.6738323333333334
$ws.Range("N10").Value = 2.021497
$ws.Range("O10").Value = 0.2263325903850996
$ws.Range("P10").Value = 0.2263325903850996
$ws.Range("Q10").Value = 2.183325696227223
$ws.Range("R10").Value = 19.649931266045
$ws.Range("S10").Value = 0.05223590254447869
$ws.Range("T10").Value = 0.05223590254447869
